$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.033147258299765
$ws.Range("D2").Value = 1.052739032815713
$ws.Range("E2").Value = 1.043560114241574
$ws.Range("F2").Value = 1.057078565411794
$ws.Range("I2").Value = 1.040685280113519
$ws.Range("J2").Value = 1.038273581367465
$ws.Range("K2").Value = 1.055486933010929
$ws.Range("L2").Value = 1.046333636256338
$ws.Range("M2").Value = 1.059814520164333
$ws.Range("N2").Value = 1.039748048944652

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.034266039114036
$ws.Range("D3").Value = 1.053493715751801
$ws.Range("E3").Value = 1.044527572754767
$ws.Range("F3").Value = 1.058070975921634
$ws.Range("I3").Value = 1.040869093023293
$ws.Range("J3").Value = 1.039034095496507
$ws.Range("K3").Value = 1.056054670911567
$ws.Range("L3").Value = 1.047111745717484
$ws.Range("M3").Value = 1.060620241870489
$ws.Range("N3").Value = 1.040509643090989

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.03498987217182
$ws.Range("D4").Value = 1.053979628154624
$ws.Range("E4").Value = 1.045153333507585
$ws.Range("F4").Value = 1.058711693124417
$ws.Range("I4").Value = 1.040985176235427
$ws.Range("J4").Value = 1.039525573601963
$ws.Range("K4").Value = 1.05641895868618
$ws.Range("L4").Value = 1.047614368637878
$ws.Range("M4").Value = 1.061139547296934
$ws.Range("N4").Value = 1.041001819151727

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.03529415012715
$ws.Range("D5").Value = 1.054183325768905
$ws.Range("E5").Value = 1.045416343450957
$ws.Range("F5").Value = 1.058980705781618
$ws.Range("I5").Value = 1.041033293142378
$ws.Range("J5").Value = 1.039732041426585
$ws.Range("K5").Value = 1.056571367709451
$ws.Range("L5").Value = 1.047825463922007
$ws.Range("M5").Value = 1.061357371938475
$ws.Range("N5").Value = 1.041208580184343

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.035345238498806
$ws.Range("D6").Value = 1.05421749345246
$ws.Range("E6").Value = 1.045460500462502
$ws.Range("F6").Value = 1.059025853955226
$ws.Range("I6").Value = 1.041041332022844
$ws.Range("J6").Value = 1.039766699515922
$ws.Range("K6").Value = 1.056596914574074
$ws.Range("L6").Value = 1.04786089558459
$ws.Range("M6").Value = 1.061393916825634
$ws.Range("N6").Value = 1.041243287492142

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.03499393803214
$ws.Range("D7").Value = 1.053982352252482
$ws.Range("E7").Value = 1.045156848094817
$ws.Range("F7").Value = 1.058715289040383
$ws.Range("I7").Value = 1.040985821866417
$ws.Range("J7").Value = 1.039528333022717
$ws.Range("K7").Value = 1.056420998082505
$ws.Range("L7").Value = 1.047617190118479
$ws.Range("M7").Value = 1.061142459811893
$ws.Range("N7").Value = 1.041004582491176

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033525375236843
$ws.Range("D8").Value = 1.052994581444529
$ws.Range("E8").Value = 1.043887122990004
$ws.Range("F8").Value = 1.057414252858239
$ws.Range("I8").Value = 1.040747991473829
$ws.Range("J8").Value = 1.038530730264364
$ws.Range("K8").Value = 1.055679439132127
$ws.Range("L8").Value = 1.046596781244452
$ws.Range("M8").Value = 1.060087241852102
$ws.Range("N8").Value = 1.040005563022479

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.030936831633346
$ws.Range("D9").Value = 1.051235530604903
$ws.Range("E9").Value = 1.041647797147215
$ws.Range("F9").Value = 1.055110657278195
$ws.Range("I9").Value = 1.040307060959437
$ws.Range("J9").Value = 1.036768028960171
$ws.Range("K9").Value = 1.054349196189753
$ws.Range("L9").Value = 1.044792060490464
$ws.Range("M9").Value = 1.058212128500543
$ws.Range("N9").Value = 1.038240358480249

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.029210577690845
$ws.Range("D10").Value = 1.05005046627587
$ws.Range("E10").Value = 1.040153630436436
$ws.Range("F10").Value = 1.053567543294338
$ws.Range("I10").Value = 1.039998459706305
$ws.Range("J10").Value = 1.035589653891015
$ws.Range("K10").Value = 1.053446604033845
$ws.Range("L10").Value = 1.043584453926705
$ws.Range("M10").Value = 1.056951527370477
$ws.Range("N10").Value = 1.037060309983334

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028462944802094
$ws.Range("D11").Value = 1.049534397860088
$ws.Range("E11").Value = 1.039506333956258
$ws.Range("F11").Value = 1.05289760771156
$ws.Range("I11").Value = 1.039861364878037
$ws.Range("J11").Value = 1.03507863125491
$ws.Range("K11").Value = 1.053052041803573
$ws.Range("L11").Value = 1.043060487434329
$ws.Range("M11").Value = 1.056403178977269
$ws.Range("N11").Value = 1.03654856163647

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.028185216300906
$ws.Range("D12").Value = 1.049342267833898
$ws.Range("E12").Value = 1.039265851999138
$ws.Range("F12").Value = 1.052648499710775
$ws.Range("I12").Value = 1.039809921118607
$ws.Range("J12").Value = 1.034888697227453
$ws.Range("K12").Value = 1.052904923394384
$ws.Range("L12").Value = 1.042865702662374
$ws.Range("M12").Value = 1.056199122393811
$ws.Range("N12").Value = 1.036358357880909

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028244791130322
$ws.Range("D13").Value = 1.049383500243694
$ws.Range("E13").Value = 1.039317438337042
$ws.Range("F13").Value = 1.052701946175819
$ws.Range("I13").Value = 1.039820979532425
$ws.Range("J13").Value = 1.034929444045145
$ws.Range("K13").Value = 1.052936506144902
$ws.Range("L13").Value = 1.042907491915566
$ws.Range("M13").Value = 1.056242910219778
$ws.Range("N13").Value = 1.036399162563756

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028439988181966
$ws.Range("D14").Value = 1.0495185252891
$ws.Range("E14").Value = 1.039486456607738
$ws.Range("F14").Value = 1.052877021754834
$ws.Range("I14").Value = 1.039857123139098
$ws.Range("J14").Value = 1.035062933648276
$ws.Range("K14").Value = 1.053039892379584
$ws.Range("L14").Value = 1.043044389735948
$ws.Range("M14").Value = 1.056386319258245
$ws.Range("N14").Value = 1.036532841737435

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.028560252261256
$ws.Range("D15").Value = 1.049601660490361
$ws.Range("E15").Value = 1.03959058807615
$ws.Range("F15").Value = 1.052984856597666
$ws.Range("I15").Value = 1.039879323431985
$ws.Range("J15").Value = 1.035145165417293
$ws.Range("K15").Value = 1.053103517804122
$ws.Range("L15").Value = 1.043128715752296
$ws.Range("M15").Value = 1.05647462853162
$ws.Range("N15").Value = 1.036615190284996

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029260192278729
$ws.Range("D16").Value = 1.050084654406964
$ws.Range("E16").Value = 1.040196582765692
$ws.Range("F16").Value = 1.053611967684298
$ws.Range("I16").Value = 1.040007485202138
$ws.Range("J16").Value = 1.035623552335818
$ws.Range("K16").Value = 1.05347271120862
$ws.Range("L16").Value = 1.043619205396196
$ws.Range("M16").Value = 1.056987866779777
$ws.Range("N16").Value = 1.037094256567817

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.029699204147376
$ws.Range("D17").Value = 1.050386840344804
$ws.Range("E17").Value = 1.040576623384316
$ws.Range("F17").Value = 1.054004867418096
$ws.Range("I17").Value = 1.040086949425411
$ws.Range("J17").Value = 1.035923423040262
$ws.Range("K17").Value = 1.053703297212126
$ws.Range("L17").Value = 1.04392659128771
$ws.Range("M17").Value = 1.057309138315106
$ws.Range("N17").Value = 1.037394553123057

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.029955257615959
$ws.Range("D18").Value = 1.050562817766142
$ws.Range("E18").Value = 1.040798264522365
$ws.Range("F18").Value = 1.054233869666412
$ws.Range("I18").Value = 1.040132964878107
$ws.Range("J18").Value = 1.036098257474533
$ws.Range("K18").Value = 1.053837433828508
$ws.Range("L18").Value = 1.044105781553534
$ws.Range("M18").Value = 1.057496289376025
$ws.Range("N18").Value = 1.037569635842278

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030042562772791
$ws.Range("D19").Value = 1.050622773578647
$ws.Range("E19").Value = 1.040873833321268
$ws.Range("F19").Value = 1.054311924764839
$ws.Range("I19").Value = 1.040148598175685
$ws.Range("J19").Value = 1.036157858779945
$ws.Range("K19").Value = 1.053883109818275
$ws.Range("L19").Value = 1.044166863388495
$ws.Range("M19").Value = 1.057560062141093
$ws.Range("N19").Value = 1.03762932178838

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029652103860938
$ws.Range("D20").Value = 1.050354447871462
$ws.Range("E20").Value = 1.040535851740373
$ws.Range("F20").Value = 1.053962730529703
$ws.Range("I20").Value = 1.040078458280768
$ws.Range("J20").Value = 1.035891257520387
$ws.Range("K20").Value = 1.053678594779297
$ws.Range("L20").Value = 1.043893622317792
$ws.Range("M20").Value = 1.057274693867028
$ws.Range("N20").Value = 1.037362341924454

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028382508198257
$ws.Range("D21").Value = 1.049478775899013
$ws.Range("E21").Value = 1.039436686245682
$ws.Range("F21").Value = 1.052825473651971
$ws.Range("I21").Value = 1.03984649411634
$ws.Range("J21").Value = 1.035023627534169
$ws.Range("K21").Value = 1.053009463176749
$ws.Range("L21").Value = 1.043004081167169
$ws.Range("M21").Value = 1.056344099238774
$ws.Range("N21").Value = 1.036493479804138

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027584121436824
$ws.Range("D22").Value = 1.048925665149497
$ws.Range("E22").Value = 1.038745323895264
$ws.Range("F22").Value = 1.052108907046987
$ws.Range("I22").Value = 1.039697637007329
$ws.Range("J22").Value = 1.034477434230477
$ws.Range("K22").Value = 1.052585511347489
$ws.Range("L22").Value = 1.042443863802981
$ws.Range("M22").Value = 1.055756824466645
$ws.Range("N22").Value = 1.035946510843305

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028007375220358
$ws.Range("D23").Value = 1.049219120287704
$ws.Range("E23").Value = 1.03911185417241
$ws.Range("F23").Value = 1.052488917456001
$ws.Range("I23").Value = 1.039776834337099
$ws.Range("J23").Value = 1.034767046257357
$ws.Range("K23").Value = 1.052810563390537
$ws.Range("L23").Value = 1.042740933653757
$ws.Range("M23").Value = 1.056068355923492
$ws.Range("N23").Value = 1.036236534152482

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029673386498025
$ws.Range("D24").Value = 1.050369085508599
$ws.Range("E24").Value = 1.040554274784085
$ws.Range("F24").Value = 1.053981770899269
$ws.Range("I24").Value = 1.040082296097558
$ws.Range("J24").Value = 1.03590579196596
$ws.Range("K24").Value = 1.053689757858448
$ws.Range("L24").Value = 1.043908519893303
$ws.Range("M24").Value = 1.057290258575415
$ws.Range("N24").Value = 1.037376897010607

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.031606126488529
$ws.Range("D25").Value = 1.051692470529431
$ws.Range("E25").Value = 1.042226942563691
$ws.Range("F25").Value = 1.055707494273608
$ws.Range("I25").Value = 1.04042363650307
$ws.Range("J25").Value = 1.037224299545594
$ws.Range("K25").Value = 1.054695877849994
$ws.Range("L25").Value = 1.045259409838406
$ws.Range("M25").Value = 1.058698746744384
$ws.Range("N25").Value = 1.038697277022238

